# New UI changes in Candidate Lobby
# - Append 4 new regression run rows (108-111) to the AMSIN sheet.
# - Fix up row 78 on the AMS sheet (style + precise run-time value).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# AMSIN sheet: append rows 108-111
# ---------------------------------------------------------------------------
$amsin = $wb.Worksheets.Item("AMSIN")

# Row 108
$amsin.Range("A108").Value = "'2023-05-09"
$amsin.Range("B108").Value = 45055.67239
$amsin.Range("B108").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$amsin.Range("C108").Value = "176htfxchanges"
$amsin.Range("D108").Value = 269
$amsin.Range("E108").Value = 267
$amsin.Range("F108").Value = 2
$amsin.Range("G108").Value = 4.38

# Row 109
$amsin.Range("A109").Value = "'2023-05-09"
$amsin.Range("B109").Value = 45055.69151628472
$amsin.Range("B109").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$amsin.Range("C109").Value = "176scndhtfx"
$amsin.Range("D109").Value = 269
$amsin.Range("E109").Value = 267
$amsin.Range("F109").Value = 2
$amsin.Range("G109").Value = 4.84

# Row 110
$amsin.Range("A110").Value = "'2023-05-11"
$amsin.Range("B110").Value = 45057.72868366898
$amsin.Range("B110").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$amsin.Range("C110").Value = "176fxhh"
$amsin.Range("D110").Value = 269
$amsin.Range("E110").Value = 265
$amsin.Range("F110").Value = 4
$amsin.Range("G110").Value = 5.17

# Row 111 - unstyled (no explicit cell style) except the Run Time cell
$amsin.Range("A111").Style = "Normal"
$amsin.Range("A111").Value = "'2023-05-12"
$amsin.Range("B111").Value = 45058.76151206988
$amsin.Range("B111").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$amsin.Range("C111").Style = "Normal"
$amsin.Range("C111").Value = "177fstcycle"
$amsin.Range("D111").Style = "Normal"
$amsin.Range("D111").Value = 269
$amsin.Range("E111").Style = "Normal"
$amsin.Range("E111").Value = 263
$amsin.Range("F111").Style = "Normal"
$amsin.Range("F111").Value = 6
$amsin.Range("G111").Style = "Normal"
$amsin.Range("G111").Value = 7.72

# ---------------------------------------------------------------------------
# AMS sheet: row 78 picks up the normal data-row style + the run time value
# is corrected to its precise timestamp.
# ---------------------------------------------------------------------------
$ams = $wb.Worksheets.Item("AMS")

$ams.Rows.Item(78).Delete()
$ams.Rows.Item(78).Insert()

$ams.Range("A78").Value = "'2023-05-08"
$ams.Range("B78").Value = 45054.53331694444
$ams.Range("B78").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ams.Range("C78").Value = "176htfxtrl"
$ams.Range("D78").Value = 269
$ams.Range("E78").Value = 266
$ams.Range("F78").Value = 3
$ams.Range("G78").Value = 4.49

Write-Output "edit complete"
